# Insert a new weekly price record for "Perejil" (Región Metropolitana /
# Mercado Mayorista Lo Valledor de Santiago) as row 634, pushing the
# existing rows 634:657 down to 635:658.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 634; this shifts rows 634-657 down to 635-658
# and carries the row-634 formatting (e.g. the date style on column D) along.
$ws.Rows.Item(634).Insert()

# Populate the new row 634 with the latest weekly observation.
$ws.Range("A634").Value = 6
$ws.Range("B634").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C634").Value = "Metropolitana"
$ws.Range("D634").Value = 44939
$ws.Range("E634").Value = 13
$ws.Range("F634").Value = 100112044
$ws.Range("G634").Value = "Perejil"
$ws.Range("H634").Value = "Sin especificar"
$ws.Range("I634").Value = "Primera"
$ws.Range("J634").Value = 330
$ws.Range("K634").Value = 10000
$ws.Range("L634").Value = 11000
$ws.Range("M634").Value = 10455
$ws.Range("N634").Value = "`$/docena de atados"
$ws.Range("O634").Value = "Región Metropolitana"
$ws.Range("P634").Value = 3485
$ws.Range("Q634").Value = 3
$ws.Range("R634").Value = "Hortaliza"
